$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Williams-Davis"
$ws.Range("A3").Value = "Gallagher-Perez"
$ws.Range("A4").Value = "Ross LLC"
$ws.Range("A5").Value = "Hull LLC"
